# Atualização de bases das ligas, do dia: 10-06-2024 às 21:53
# Swap the data (columns B:AD) between rows 26/27 and between rows 160/161,
# leaving column A (the running id/index) untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Swap-RowData {
    param($Worksheet, $Row1, $Row2)

    $rng1 = $Worksheet.Range("B$Row1`:AD$Row1")
    $rng2 = $Worksheet.Range("B$Row2`:AD$Row2")

    $vals1 = $rng1.Value2
    $vals2 = $rng2.Value2

    $rng1.Value2 = $vals2
    $rng2.Value2 = $vals1
}

Swap-RowData $ws 26 27
Swap-RowData $ws 160 161
